# Helper: write a value to a cell as TEXT (shared string), matching the
# source workbook's convention of storing numeric-looking labels as text,
# without leaving a lingering custom number format on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# --- Sheet OM ---
$ws = $wb.Worksheets.Item("OM")
$ws.Range("B2").Value = 78.85
$ws.Range("B3").Value = 19.19
$ws.Range("B4").Value = 1.54
$ws.Range("B5").Value = 0.28
Set-TextValue $ws.Range("A6") "4"
$ws.Range("B6").Value = 0.14

# --- Sheet NV ---
$ws = $wb.Worksheets.Item("NV")
$ws.Range("B2").Value = 94.68
$ws.Range("B3").Value = 5.18

# --- Sheet NR ---
$ws = $wb.Worksheets.Item("NR")
$ws.Range("B2").Value = 70.73
$ws.Range("B3").Value = 12.32
$ws.Range("B4").Value = 7.28
$ws.Range("B5").Value = 4.62
$ws.Range("B6").Value = 2.94
$ws.Range("B7").Value = 0.56
$ws.Range("B8").Value = 1.26
$ws.Range("B9").Value = 0.28

# --- Sheet ALL ---
$ws = $wb.Worksheets.Item("ALL")
$ws.Range("B2").Value = 50.7
$ws.Range("B3").Value = 27.73
$ws.Range("B4").Value = 8.96
$ws.Range("B5").Value = 6.44
$ws.Range("B6").Value = 3.78
$ws.Range("B7").Value = 0.7
$ws.Range("B8").Value = 1.4
$ws.Range("B9").Value = 0.28
# Remove row 10 (was A10="8", B10=0.14) entirely
$ws.Range("A10:B10").Delete()

# --- Sheet summary ---
$ws = $wb.Worksheets.Item("summary")
Set-TextValue $ws.Range("A2") "0.24"
Set-TextValue $ws.Range("B2") "0.05"
Set-TextValue $ws.Range("C2") "0.65"
Set-TextValue $ws.Range("E2") "0.94"

Set-TextValue $ws.Range("A3") "0.5"
Set-TextValue $ws.Range("B3") "0.23"
Set-TextValue $ws.Range("C3") "1.27"
Set-TextValue $ws.Range("E3") "1.33"

Set-TextValue $ws.Range("A6") "4"
Set-TextValue $ws.Range("E6") "7"
